# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> currently "Office Theme" (used only by the notes master)
#   ppt/theme/theme2.xml -> currently "Integral"/"Red Violet" (used by the slide master,
#                            and this is the one PowerPoint's object model exposes as
#                            ActivePresentation.SlideMaster.Theme)
#
# The target edit swaps the two themes' contents: the slide-master-facing theme part
# should end up holding the "Office Theme" colour scheme (what theme1.xml has today),
# while the notes-master-facing theme part should end up holding the "Integral"/
# "Red Violet" colours (what theme2.xml has today). Font scheme and format scheme are
# identical between the two themes already, so only the 12-colour scheme differs.
#
# Recolour the slide master's theme colour scheme to the "Office Theme" palette so the
# reachable theme part (theme2.xml) matches the post-swap target content exactly.
# ColorFormat.RGB uses the standard OLE (BGR-packed) colour order, so each target
# sRGB hex value below is converted accordingly.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0        # dk1      sRGB 000000
$cs.Item(2).RGB  = 16777215 # lt1      sRGB FFFFFF
$cs.Item(3).RGB  = 6968388  # dk2      sRGB 44546A
$cs.Item(4).RGB  = 15132391 # lt2      sRGB E7E6E6
$cs.Item(5).RGB  = 13998939 # accent1  sRGB 5B9BD5
$cs.Item(6).RGB  = 3243501  # accent2  sRGB ED7D31
$cs.Item(7).RGB  = 10855845 # accent3  sRGB A5A5A5
$cs.Item(8).RGB  = 49407    # accent4  sRGB FFC000
$cs.Item(9).RGB  = 12874308 # accent5  sRGB 4472C4
$cs.Item(10).RGB = 4697456  # accent6  sRGB 70AD47
$cs.Item(11).RGB = 12673797 # hlink    sRGB 0563C1
$cs.Item(12).RGB = 7491477  # folHlink sRGB 954F72
